# "finally fixed the actuals"
# Adds a new "KPIType" column (Transaction / Snapshot) to the d_kpi sheet,
# fixes AggregateNum for the Totalvaluelocked KPI (row 11) and turns off the
# "live" flag for the last KPI (row 17, now column X after the insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("d_kpi")

# 1) Insert a new column before the old "Level0_Attribuut" column (O).
#    Everything from O onward shifts one column to the right (O -> P, ... X -> Y).
$ws.Range("O1").EntireColumn.Insert()

# Give the freshly inserted column roughly the same display width as its
# neighbours (customWidth, no bestFit) - closest width this host can store.
$ws.Columns.Item(15).ColumnWidth = 14.42

# 2) Fill in the new "KPIType" column.
$ws.Range("O1").Value = "KPIType"

$lastRow = 17
for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 11) {
        $ws.Cells.Item($r, 15).Value = "Snapshot"
    } else {
        $ws.Cells.Item($r, 15).Value = "Transaction"
    }
}

# 3) Fix AggregateNum (column M) for KPI #10 (Totalvaluelocked), row 11: 2 -> 1.
$ws.Range("M11").Value = 1

# 4) Turn off "live" (now column X after the insert) for KPI #16, row 17: 1 -> 0.
$ws.Range("X17").Value = 0

# 5) Update the view a little (best effort - scroll right, select last cell).
$ws.Range("X17").Select()

$wb.Save()
